$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 48461
$ws.Cells.Item(2, 5).Value = 951190896285
$ws.Cells.Item(2, 6).Value = 19786255673
$ws.Cells.Item(2, 7).Value = 2.52235

# Row 3
$ws.Cells.Item(3, 4).Value = 2531
$ws.Cells.Item(3, 5).Value = 304242038355
$ws.Cells.Item(3, 6).Value = 7617326047
$ws.Cells.Item(3, 7).Value = 1.61553

# Row 4
$ws.Cells.Item(4, 4).Value = 1.001
$ws.Cells.Item(4, 5).Value = 96454290308
$ws.Cells.Item(4, 6).Value = 26084958286
$ws.Cells.Item(4, 7).Value = -0.01342

# Row 5
$ws.Cells.Item(5, 4).Value = 322.92
$ws.Cells.Item(5, 5).Value = 49701832805
$ws.Cells.Item(5, 6).Value = 476923714
$ws.Cells.Item(5, 7).Value = 0.56041

# Row 6
$ws.Cells.Item(6, 4).Value = 110.18
$ws.Cells.Item(6, 5).Value = 48126993126
$ws.Cells.Item(6, 6).Value = 1469933219
$ws.Cells.Item(6, 7).Value = 2.01941

# Row 7
$ws.Cells.Item(7, 4).Value = 0.533124
$ws.Cells.Item(7, 5).Value = 29068901850
$ws.Cells.Item(7, 6).Value = 620579571
$ws.Cells.Item(7, 7).Value = 2.16463

# Row 8
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 27968836234
$ws.Cells.Item(8, 6).Value = 4219827484
$ws.Cells.Item(8, 7).Value = 0.01127

# Row 9
$ws.Cells.Item(9, 4).Value = 2530.06
$ws.Cells.Item(9, 5).Value = 24350888654
$ws.Cells.Item(9, 6).Value = 10464147
$ws.Cells.Item(9, 7).Value = 1.72665

# Row 10
$ws.Cells.Item(10, 4).Value = 0.552951
$ws.Cells.Item(10, 5).Value = 19426771679
$ws.Cells.Item(10, 6).Value = 517437611
$ws.Cells.Item(10, 7).Value = 3.43052

# Row 11
$ws.Cells.Item(11, 4).Value = 40.62
$ws.Cells.Item(11, 5).Value = 14911805042
$ws.Cells.Item(11, 6).Value = 643057134
$ws.Cells.Item(11, 7).Value = 5.10293

# Row 12
$ws.Cells.Item(12, 2).Value = 'LINK'
$ws.Cells.Item(12, 3).Value = 'Chainlink'
$ws.Cells.Item(12, 4).Value = 20.47
$ws.Cells.Item(12, 5).Value = 12022656082
$ws.Cells.Item(12, 6).Value = 921624967
$ws.Cells.Item(12, 7).Value = 12.23084

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'DOGE'
$ws.Cells.Item(13, 3).Value = 'Dogecoin'
$ws.Cells.Item(13, 4).Value = 0.082359
$ws.Cells.Item(13, 5).Value = 11780265319
$ws.Cells.Item(13, 6).Value = 330291370
$ws.Cells.Item(13, 7).Value = 1.79812

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 4).Value = 0.1248
$ws.Cells.Item(14, 5).Value = 10994834803
$ws.Cells.Item(14, 6).Value = 199427326
$ws.Cells.Item(14, 7).Value = 1.15001

# Row 15
$ws.Cells.Item(15, 4).Value = 7.3
$ws.Cells.Item(15, 5).Value = 9698771289
$ws.Cells.Item(15, 6).Value = 167873726
$ws.Cells.Item(15, 7).Value = 2.46307

# Row 16
$ws.Cells.Item(16, 4).Value = 0.856474
$ws.Cells.Item(16, 5).Value = 7951193099
$ws.Cells.Item(16, 6).Value = 261032112
$ws.Cells.Item(16, 7).Value = 1.25351

# Row 17
$ws.Cells.Item(17, 2).Value = 'WBTC'
$ws.Cells.Item(17, 3).Value = 'Wrapped Bitcoin'
$ws.Cells.Item(17, 4).Value = 48277
$ws.Cells.Item(17, 5).Value = 7620494088
$ws.Cells.Item(17, 6).Value = 208042778
$ws.Cells.Item(17, 7).Value = 2.26768

# Row 18
$ws.Cells.Item(18, 2).Value = 'TON'
$ws.Cells.Item(18, 3).Value = 'Toncoin'
$ws.Cells.Item(18, 4).Value = 2.12
$ws.Cells.Item(18, 5).Value = 7331859545
$ws.Cells.Item(18, 6).Value = 15461018
$ws.Cells.Item(18, 7).Value = 2.31448

# Row 19
$ws.Cells.Item(19, 4).Value = 13.48
$ws.Cells.Item(19, 5).Value = 6181614896
$ws.Cells.Item(19, 6).Value = 112206421
$ws.Cells.Item(19, 7).Value = 5.70007

# Row 20
$ws.Cells.Item(20, 4).Value = 0.0000095
$ws.Cells.Item(20, 5).Value = 5599871555
$ws.Cells.Item(20, 6).Value = 111815428
$ws.Cells.Item(20, 7).Value = 1.89579

# Row 21
$ws.Cells.Item(21, 4).Value = 71.98999999999999
$ws.Cells.Item(21, 5).Value = 5343077689
$ws.Cells.Item(21, 6).Value = 286191676
$ws.Cells.Item(21, 7).Value = 2.25266

# Row 22
$ws.Cells.Item(22, 2).Value = 'BCH'
$ws.Cells.Item(22, 3).Value = 'Bitcoin Cash'
$ws.Cells.Item(22, 4).Value = 266.82
$ws.Cells.Item(22, 5).Value = 5216378852
$ws.Cells.Item(22, 6).Value = 408583575
$ws.Cells.Item(22, 7).Value = 8.665179999999999

# Row 23
$ws.Cells.Item(23, 2).Value = 'UNI'
$ws.Cells.Item(23, 3).Value = 'Uniswap'
$ws.Cells.Item(23, 4).Value = 6.64
$ws.Cells.Item(23, 5).Value = 5005779821
$ws.Cells.Item(23, 6).Value = 60627289
$ws.Cells.Item(23, 7).Value = 0.43463

# Row 24
$ws.Cells.Item(24, 2).Value = 'DAI'
$ws.Cells.Item(24, 3).Value = 'Dai'
$ws.Cells.Item(24, 4).Value = 0.998206
$ws.Cells.Item(24, 5).Value = 4917520805
$ws.Cells.Item(24, 6).Value = 228913794
$ws.Cells.Item(24, 7).Value = -0.23665

# Row 25
$ws.Cells.Item(25, 2).Value = 'ATOM'
$ws.Cells.Item(25, 3).Value = 'Cosmos Hub'
$ws.Cells.Item(25, 4).Value = 10.16
$ws.Cells.Item(25, 5).Value = 3900295927
$ws.Cells.Item(25, 6).Value = 170728646
$ws.Cells.Item(25, 7).Value = 1.45881

# Row 26
$ws.Cells.Item(26, 2).Value = 'LEO'
$ws.Cells.Item(26, 3).Value = 'LEO Token'
$ws.Cells.Item(26, 4).Value = 4.12
$ws.Cells.Item(26, 5).Value = 3818813705
$ws.Cells.Item(26, 6).Value = 1437211
$ws.Cells.Item(26, 7).Value = 0.43726

# Row 27
$ws.Cells.Item(27, 2).Value = 'ETC'
$ws.Cells.Item(27, 3).Value = 'Ethereum Classic'
$ws.Cells.Item(27, 4).Value = 26.08
$ws.Cells.Item(27, 5).Value = 3735683370
$ws.Cells.Item(27, 6).Value = 130285848
$ws.Cells.Item(27, 7).Value = 1.49157

# Row 28
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'IMX'
$ws.Cells.Item(28, 3).Value = 'Immutable'
$ws.Cells.Item(28, 4).Value = 2.7
$ws.Cells.Item(28, 5).Value = 3658965022
$ws.Cells.Item(28, 6).Value = 148822987
$ws.Cells.Item(28, 7).Value = -1.18029

# Row 29
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 'TAO'
$ws.Cells.Item(29, 3).Value = 'Bittensor'
$ws.Cells.Item(29, 4).Value = 569.52
$ws.Cells.Item(29, 5).Value = 3554949878
$ws.Cells.Item(29, 6).Value = 41851065
$ws.Cells.Item(29, 7).Value = 3.33664

# Row 30
$ws.Cells.Item(30, 2).Value = 'OP'
$ws.Cells.Item(30, 3).Value = 'Optimism'
$ws.Cells.Item(30, 4).Value = 3.63
$ws.Cells.Item(30, 5).Value = 3476114169
$ws.Cells.Item(30, 6).Value = 172037473
$ws.Cells.Item(30, 7).Value = 7.54606

# Row 31
$ws.Cells.Item(31, 2).Value = 'NEAR'
$ws.Cells.Item(31, 3).Value = 'NEAR Protocol'
$ws.Cells.Item(31, 4).Value = 3.16
$ws.Cells.Item(31, 5).Value = 3275480918
$ws.Cells.Item(31, 6).Value = 190974110
$ws.Cells.Item(31, 7).Value = 5.03088

# Row 32
$ws.Cells.Item(32, 2).Value = 'KAS'
$ws.Cells.Item(32, 3).Value = 'Kaspa'
$ws.Cells.Item(32, 4).Value = 0.144584
$ws.Cells.Item(32, 5).Value = 3275414704
$ws.Cells.Item(32, 6).Value = 50508351
$ws.Cells.Item(32, 7).Value = 6.14877

# Row 33
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 'TIA'
$ws.Cells.Item(33, 3).Value = 'Celestia'
$ws.Cells.Item(33, 4).Value = 19.75
$ws.Cells.Item(33, 5).Value = 3236508845
$ws.Cells.Item(33, 6).Value = 152970097
$ws.Cells.Item(33, 7).Value = -2.39467

# Row 34
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 'XLM'
$ws.Cells.Item(34, 3).Value = 'Stellar'
$ws.Cells.Item(34, 4).Value = 0.112147
$ws.Cells.Item(34, 5).Value = 3186080446
$ws.Cells.Item(34, 6).Value = 53067855
$ws.Cells.Item(34, 7).Value = 0.7483300000000001

# Row 35
$ws.Cells.Item(35, 2).Value = 'INJ'
$ws.Cells.Item(35, 3).Value = 'Injective'
$ws.Cells.Item(35, 4).Value = 35.8
$ws.Cells.Item(35, 5).Value = 3168916238
$ws.Cells.Item(35, 6).Value = 156211623
$ws.Cells.Item(35, 7).Value = 4.13752

# Row 36
$ws.Cells.Item(36, 2).Value = 'APT'
$ws.Cells.Item(36, 3).Value = 'Aptos'
$ws.Cells.Item(36, 4).Value = 9.18
$ws.Cells.Item(36, 5).Value = 3117772396
$ws.Cells.Item(36, 6).Value = 78437410
$ws.Cells.Item(36, 7).Value = 0.92775

# Row 37
$ws.Cells.Item(37, 2).Value = 'OKB'
$ws.Cells.Item(37, 3).Value = 'OKB'
$ws.Cells.Item(37, 4).Value = 49.75
$ws.Cells.Item(37, 5).Value = 2985974343
$ws.Cells.Item(37, 6).Value = 5012442
$ws.Cells.Item(37, 7).Value = 0.37987

# Row 38
$ws.Cells.Item(38, 4).Value = 1.005
$ws.Cells.Item(38, 5).Value = 2773842649
$ws.Cells.Item(38, 6).Value = 3970353114
$ws.Cells.Item(38, 7).Value = -0.01832

# Row 39
$ws.Cells.Item(39, 2).Value = 'STX'
$ws.Cells.Item(39, 3).Value = 'Stacks'
$ws.Cells.Item(39, 4).Value = 1.91
$ws.Cells.Item(39, 5).Value = 2746462690
$ws.Cells.Item(39, 6).Value = 81281031
$ws.Cells.Item(39, 7).Value = 8.250690000000001

# Row 40
$ws.Cells.Item(40, 2).Value = 'FIL'
$ws.Cells.Item(40, 3).Value = 'Filecoin'
$ws.Cells.Item(40, 4).Value = 5.41
$ws.Cells.Item(40, 5).Value = 2741216455
$ws.Cells.Item(40, 6).Value = 109682438
$ws.Cells.Item(40, 7).Value = 1.01905

# Row 41
$ws.Cells.Item(41, 2).Value = 'LDO'
$ws.Cells.Item(41, 3).Value = 'Lido DAO'
$ws.Cells.Item(41, 4).Value = 3.03
$ws.Cells.Item(41, 5).Value = 2697684726
$ws.Cells.Item(41, 6).Value = 52531883
$ws.Cells.Item(41, 7).Value = 3.14664

# Row 42
$ws.Cells.Item(42, 2).Value = 'HBAR'
$ws.Cells.Item(42, 3).Value = 'Hedera'
$ws.Cells.Item(42, 4).Value = 0.07879
$ws.Cells.Item(42, 5).Value = 2654078198
$ws.Cells.Item(42, 6).Value = 29248572
$ws.Cells.Item(42, 7).Value = 0.9248499999999999

# Row 43
$ws.Cells.Item(43, 4).Value = 2
$ws.Cells.Item(43, 5).Value = 2549867035
$ws.Cells.Item(43, 6).Value = 193200366
$ws.Cells.Item(43, 7).Value = 1.97692

# Row 44
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 4).Value = 0.088699
$ws.Cells.Item(44, 5).Value = 2354128890
$ws.Cells.Item(44, 6).Value = 6519621
$ws.Cells.Item(44, 7).Value = 1.68057

# Row 45
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'XMR'
$ws.Cells.Item(45, 3).Value = 'Monero'
$ws.Cells.Item(45, 4).Value = 121.47
$ws.Cells.Item(45, 5).Value = 2202964980
$ws.Cells.Item(45, 6).Value = 79035116
$ws.Cells.Item(45, 7).Value = 1.89893

# Row 46
$ws.Cells.Item(46, 4).Value = 0.0301961
$ws.Cells.Item(46, 5).Value = 2196189686
$ws.Cells.Item(46, 6).Value = 36479329
$ws.Cells.Item(46, 7).Value = 2.37214

# Row 47
$ws.Cells.Item(47, 4).Value = 0.66534
$ws.Cells.Item(47, 5).Value = 2108005683
$ws.Cells.Item(47, 6).Value = 37928721
$ws.Cells.Item(47, 7).Value = 3.48303

# Row 48
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 'SUI'
$ws.Cells.Item(48, 3).Value = 'Sui'
$ws.Cells.Item(48, 4).Value = 1.73
$ws.Cells.Item(48, 5).Value = 2012732904
$ws.Cells.Item(48, 6).Value = 248170215
$ws.Cells.Item(48, 7).Value = 0.7506699999999999

# Row 49
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 4).Value = 2029.85
$ws.Cells.Item(49, 5).Value = 1872841087
$ws.Cells.Item(49, 6).Value = 47817925
$ws.Cells.Item(49, 7).Value = 2.23151

# Row 50
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = 'RNDR'
$ws.Cells.Item(50, 3).Value = 'Render'
$ws.Cells.Item(50, 4).Value = 4.71
$ws.Cells.Item(50, 5).Value = 1772528355
$ws.Cells.Item(50, 6).Value = 74260857
$ws.Cells.Item(50, 7).Value = 1.90401

# Row 51
$ws.Cells.Item(51, 4).Value = 0.690954
$ws.Cells.Item(51, 5).Value = 1675176226
$ws.Cells.Item(51, 6).Value = 110078019
$ws.Cells.Item(51, 7).Value = 0.88615
